# MAI_holdings.xlsx update:
#   - Roll the "as of" date in the confidential disclosure footer forward
#     from 2021-03-31 to 2021-04-05.
#   - Refresh the Weight / Percent Change figures for the five holdings
#     rows plus the Total row's Percent Change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected (accidental-edit guard only); unprotect so the
# cell writes below are allowed, then restore protection afterwards.
$ws.Unprotect()

$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-05 for illustrative purposes only and are subject to change."
# Re-fit row 10 so the wrapped two-line disclosure text doesn't leave a
# stale/custom row height behind.
$ws.Rows(10).EntireRow.AutoFit()

$ws.Range("D2").Value = 0.489047130225389
$ws.Range("E2").Value = 0.001956181533646406

$ws.Range("D3").Value = 0.3325166015987633
$ws.Range("E3").Value = 0.007381889763779625

$ws.Range("D4").Value = 0.09424660431141206
$ws.Range("E4").Value = -0.0001869857890800519

$ws.Range("D5").Value = 0.05480053850167384
$ws.Range("E5").Value = 0.0004594004823705866

$ws.Range("D6").Value = 0.02938912536276169
$ws.Range("E6").Value = -0.003854802441374816

$ws.Range("E7").Value = 0.003305529208806446

$ws.Protect()
